$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 17424.8
$ws.Range("I20").Value = 4700
$ws.Range("J20").Value = 36512
$ws.Range("K20").Value = 4700
$ws.Range("L20").Value = 36512
$ws.Range("M20").Value = -4470
$ws.Range("N20").Value = -36972
$ws.Range("H29").Value = 952.5
$ws.Range("I29").Value = 176.42857
$ws.Range("J29").Value = 1728.5714
$ws.Range("K29").Value = 529.28571
$ws.Range("L29").Value = 5185.7142
$ws.Range("M29").Value = -248.28571
$ws.Range("N29").Value = -5747.7142
$ws.Range("H35").Value = 17424.8
$ws.Range("I35").Value = 4700
$ws.Range("J35").Value = 36512
$ws.Range("K35").Value = 4700
$ws.Range("L35").Value = 36512
$ws.Range("M35").Value = -4321
$ws.Range("N35").Value = -37270
$ws.Range("H58").Value = 3404
$ws.Range("I58").Value = 447
$ws.Range("J58").Value = 5252.125
$ws.Range("K58").Value = 1341
$ws.Range("L58").Value = 15756.375
$ws.Range("M58").Value = -1191
$ws.Range("N58").Value = -16056.375
$ws.Range("H81").Value = 40000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 40000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H84").Value = 40000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 40000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H107").Value = 141.63637
$ws.Range("I107").Value = 141.63637
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 141.63637
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1778.36363
$ws.Range("H116").Value = 3282.3
$ws.Range("I116").Value = 3395.3845
$ws.Range("J116").Value = 3072.2856
$ws.Range("K116").Value = 3395.3845
$ws.Range("L116").Value = 3072.2856
$ws.Range("M116").Value = 46.61549999999988
$ws.Range("N116").Value = -9956.285599999999
$ws.Range("H138").Value = 2686.6768
$ws.Range("I138").Value = 643.2759
$ws.Range("J138").Value = 3533.2285
$ws.Range("K138").Value = 1929.8277
$ws.Range("L138").Value = 10599.6855
$ws.Range("M138").Value = 3210.1723
$ws.Range("N138").Value = -20879.6855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2849.5
$ws.Range("I132").Value = 1833.3334
$ws.Range("J132").Value = 3285
$ws.Range("K132").Value = 5500.0002
$ws.Range("L132").Value = 9855
$ws.Range("M132").Value = -2970.0002
$ws.Range("N132").Value = -14915

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 29200
$ws.Range("H134").Value = 27216.666
$ws.Range("I134").Value = 3861.4092
$ws.Range("J134").Value = 52907.45
$ws.Range("K134").Value = 11584.2276
$ws.Range("L134").Value = 158722.35
$ws.Range("M134").Value = -9049.2276
$ws.Range("N134").Value = -163792.35

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 903.5
$ws.Range("I10").Value = 903.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 903.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -764.5
$ws.Range("H62").Value = 4347.353
$ws.Range("I62").Value = 2772.1428
$ws.Range("J62").Value = 5450
$ws.Range("K62").Value = 2772.1428
$ws.Range("L62").Value = 5450
$ws.Range("M62").Value = -2148.1428
$ws.Range("N62").Value = -6698
$ws.Range("H65").Value = 4347.353
$ws.Range("I65").Value = 2772.1428
$ws.Range("J65").Value = 5450
$ws.Range("K65").Value = 13860.714
$ws.Range("L65").Value = 27250
$ws.Range("M65").Value = -10740.714
$ws.Range("N65").Value = -33490
$ws.Range("H81").Value = 33900
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 33900
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 33900
$ws.Range("N81").Value = -35896
$ws.Range("H84").Value = 33900
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 33900
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 101700
$ws.Range("N84").Value = -111684
$ws.Range("H140").Value = 53558.184
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 53558.184
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 53558.184
$ws.Range("N140").Value = -63918.184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 681.25
$ws.Range("I5").Value = 672.2727
$ws.Range("J5").Value = 780
$ws.Range("K5").Value = 2016.8181
$ws.Range("L5").Value = 2340
$ws.Range("M5").Value = -1904.8181
$ws.Range("N5").Value = -2564
$ws.Range("H86").Value = 1340
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1340
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4020
$ws.Range("N86").Value = -6392
$ws.Range("H89").Value = 1340
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1340
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 12060
$ws.Range("N89").Value = -23916
$ws.Range("H117").Value = 2313.8333
$ws.Range("I117").Value = 265
$ws.Range("J117").Value = 2723.6
$ws.Range("K117").Value = 795
$ws.Range("L117").Value = 8170.799999999999
$ws.Range("M117").Value = 2647
$ws.Range("N117").Value = -15054.8
$ws.Range("H135").Value = 681.25
$ws.Range("I135").Value = 672.2727
$ws.Range("J135").Value = 780
$ws.Range("K135").Value = 6050.454299999999
$ws.Range("L135").Value = 7020
$ws.Range("M135").Value = -3515.454299999999
$ws.Range("N135").Value = -12090

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 1690934.4
$ws.Range("I18").Value = 3338666.8
$ws.Range("J18").Value = 43202
$ws.Range("K18").Value = 3338666.8
$ws.Range("L18").Value = 43202
$ws.Range("M18").Value = -3338373.8
$ws.Range("N18").Value = -43788
$ws.Range("H43").Value = 4500
$ws.Range("I43").Value = 4500
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 4500
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -4349
$ws.Range("N43").ClearContents()
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 246
$ws.Range("N99").ClearContents()
$ws.Range("H132").Value = 3372.8572
$ws.Range("I132").Value = 856
$ws.Range("J132").Value = 4379.6
$ws.Range("K132").Value = 2568
$ws.Range("L132").Value = 13138.8
$ws.Range("M132").Value = -38
$ws.Range("N132").Value = -18198.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1858.3636
$ws.Range("I7").Value = 1798.5
$ws.Range("J7").Value = 1963.125
$ws.Range("K7").Value = 1798.5
$ws.Range("L7").Value = 1963.125
$ws.Range("M7").Value = -1686.5
$ws.Range("N7").Value = -2187.125
$ws.Range("H126").Value = 1858.3636
$ws.Range("I126").Value = 1798.5
$ws.Range("J126").Value = 1963.125
$ws.Range("K126").Value = 5395.5
$ws.Range("L126").Value = 5889.375
$ws.Range("M126").Value = -2925.5
$ws.Range("N126").Value = -10829.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 9500
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 9500
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 9500
$ws.Range("N47").Value = -10644
$ws.Range("H107").Value = 571.8125
$ws.Range("I107").Value = 620.4
$ws.Range("J107").Value = 398.2857
$ws.Range("K107").Value = 1861.2
$ws.Range("L107").Value = 1194.8571
$ws.Range("M107").Value = 58.80000000000018
$ws.Range("H136").Value = 1662.303
$ws.Range("I136").Value = 1560.3793
$ws.Range("J136").Value = 2401.25
$ws.Range("K136").Value = 4681.1379
$ws.Range("L136").Value = 7203.75
$ws.Range("M136").Value = -2131.1379
